$d = $word.ActiveDocument

$replacements = @(
    @("30÷7=", "28÷8="),
    @("42÷7=", "64÷6="),
    @("45÷9=", "77÷8="),
    @("84÷3=", "94÷2="),
    @("82÷8=", "63÷4="),
    @("35÷4=", "55÷4="),
    @("94÷6=", "96÷6="),
    @("25÷3=", "33÷2="),
    @("36÷8=", "16÷3="),
    @("19÷4=", "86÷6="),
    @("65÷6=", "19÷5="),
    @("21÷2=", "85÷8="),
    @("66÷4=", "92÷9="),
    @("98÷5=", "46÷4="),
    @("76÷6=", "20÷8="),
    @("81÷9=", "52÷8="),
    @("24÷8=", "19÷4="),
    @("65÷4=", "34÷5="),
    @("13÷9=", "36÷6="),
    @("43÷8=", "30÷8="),
    @("28÷3=", "49÷2="),
    @("74÷3=", "71÷9="),
    @("57÷5=", "91÷4="),
    @("41÷7=", "47÷5="),
    @("13÷5=", "17÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
